$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 57124.332
$ws.Range("J21").Value = 34229.668
$ws.Range("L21").Value = 34229.668
$ws.Range("N21").Value = -35165.668

$ws.Range("H23").Value = 57124.332
$ws.Range("J23").Value = 34229.668
$ws.Range("L23").Value = 34229.668
$ws.Range("N23").Value = -34697.668

$ws.Range("H51").Value = 7840
$ws.Range("J51").Value = 7840
$ws.Range("L51").Value = 7840
$ws.Range("N51").Value = -8808

$ws.Range("H116").Value = 13109.308
$ws.Range("I116").Value = 3538.3333
$ws.Range("J116").Value = 15980.6
$ws.Range("K116").Value = 3538.3333
$ws.Range("L116").Value = 15980.6
$ws.Range("M116").Value = -96.33329999999978
$ws.Range("N116").Value = -22864.6

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H129").Value = 1017.13043
$ws.Range("I129").Value = 333.33334
$ws.Range("J129").Value = 1048.2122
$ws.Range("K129").Value = 1000.00002
$ws.Range("L129").Value = 3144.6366
$ws.Range("M129").Value = 3999.99998
$ws.Range("N129").Value = -13144.6366

$ws.Range("H132").Value = 73157.836
$ws.Range("I132").Value = 86215.66
$ws.Range("K132").Value = 258646.98
$ws.Range("M132").Value = -256116.98

$ws.Range("H137").Value = 3121.5557
$ws.Range("I137").Value = 2006.0588
$ws.Range("J137").Value = 5017.9
$ws.Range("K137").Value = 6018.1764
$ws.Range("L137").Value = 15053.7
$ws.Range("M137").Value = -3468.1764
$ws.Range("N137").Value = -20153.7

$ws.Range("H138").Value = 3377.602
$ws.Range("I138").Value = 1273.3182
$ws.Range("J138").Value = 3986.7368
$ws.Range("K138").Value = 3819.9546
$ws.Range("L138").Value = 11960.2104
$ws.Range("M138").Value = 1320.0454
$ws.Range("N138").Value = -22240.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7923.2
$ws.Range("I32").Value = 5601.7144
$ws.Range("J32").Value = 10489.053
$ws.Range("K32").Value = 5601.7144
$ws.Range("L32").Value = 10489.053
$ws.Range("M32").Value = -5314.7144
$ws.Range("N32").Value = -11063.053

$ws.Range("H45").Value = 976.86365
$ws.Range("I45").Value = 844.63635
$ws.Range("J45").Value = 1109.091
$ws.Range("K45").Value = 844.63635
$ws.Range("L45").Value = 1109.091
$ws.Range("M45").Value = -467.63635
$ws.Range("N45").Value = -1863.091

$ws.Range("H60").Value = 26025.5
$ws.Range("I60").Value = 2051
$ws.Range("J60").Value = 50000
$ws.Range("K60").Value = 2051
$ws.Range("L60").Value = 50000
$ws.Range("M60").Value = -1318
$ws.Range("N60").Value = -51466

$ws.Range("H63").Value = 9237714
$ws.Range("I63").Value = 23087436
$ws.Range("J63").Value = 4566.6665
$ws.Range("K63").Value = 23087436
$ws.Range("L63").Value = 4566.6665
$ws.Range("M63").Value = -23086750
$ws.Range("N63").Value = -5938.6665

$ws.Range("H66").Value = 9237714
$ws.Range("I66").Value = 23087436
$ws.Range("J66").Value = 4566.6665
$ws.Range("K66").Value = 115437180
$ws.Range("L66").Value = 22833.3325
$ws.Range("M66").Value = -115433748
$ws.Range("N66").Value = -29697.3325

$ws.Range("H68").Value = 49999.5
$ws.Range("J68").Value = 49999.5
$ws.Range("L68").Value = 49999.5
$ws.Range("N68").Value = -51621.5

$ws.Range("H71").Value = 49999.5
$ws.Range("J71").Value = 49999.5
$ws.Range("L71").Value = 149998.5
$ws.Range("N71").Value = -158110.5

$ws.Range("H74").Value = 1418.5333
$ws.Range("I74").Value = 1066.909
$ws.Range("J74").Value = 2385.5
$ws.Range("K74").Value = 1066.909
$ws.Range("L74").Value = 2385.5
$ws.Range("M74").Value = -192.9090000000001
$ws.Range("N74").Value = -4133.5

$ws.Range("H77").Value = 1418.5333
$ws.Range("I77").Value = 1066.909
$ws.Range("J77").Value = 2385.5
$ws.Range("K77").Value = 5334.545
$ws.Range("L77").Value = 11927.5
$ws.Range("M77").Value = -966.5450000000001
$ws.Range("N77").Value = -20663.5

$ws.Range("H97").Value = 763.38464
$ws.Range("I97").Value = 602.087
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 602.087
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -106.087
$ws.Range("N97").Value = -2992

$ws.Range("H132").Value = 3235.0667
$ws.Range("I132").Value = 2166.158
$ws.Range("J132").Value = 5081.364
$ws.Range("K132").Value = 6498.474
$ws.Range("L132").Value = 15244.092
$ws.Range("M132").Value = -3968.474
$ws.Range("N132").Value = -20304.092

$ws.Range("H137").Value = 40497.5
$ws.Range("J137").Value = 40497.5
$ws.Range("L137").Value = 40497.5
$ws.Range("N137").Value = -50697.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 40000
$ws.Range("J19").Value = 40000
$ws.Range("L19").Value = 40000
$ws.Range("N19").Value = -40346

$ws.Range("H59").Value = 118874
$ws.Range("J59").Value = 118874
$ws.Range("L59").Value = 118874
$ws.Range("N59").Value = -120568

$ws.Range("H64").Value = 367.22223
$ws.Range("I64").Value = 226.5
$ws.Range("J64").Value = 479.8
$ws.Range("K64").Value = 226.5
$ws.Range("L64").Value = 479.8
$ws.Range("M64").Value = -1.5
$ws.Range("N64").Value = -929.8

$ws.Range("H67").Value = 367.22223
$ws.Range("I67").Value = 226.5
$ws.Range("J67").Value = 479.8
$ws.Range("K67").Value = 226.5
$ws.Range("L67").Value = 479.8
$ws.Range("M67").Value = 553.5
$ws.Range("N67").Value = -2039.8

$ws.Range("H94").Value = 1330.1428
$ws.Range("I94").Value = 978.13635
$ws.Range("J94").Value = 2620.8333
$ws.Range("K94").Value = 978.13635
$ws.Range("L94").Value = 2620.8333
$ws.Range("M94").Value = -527.13635
$ws.Range("N94").Value = -3522.8333

$ws.Range("H137").Value = 33000
$ws.Range("J137").Value = 41000
$ws.Range("L137").Value = 41000
$ws.Range("N137").Value = -51200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17245344
$ws.Range("I31").Value = 1986.7333
$ws.Range("J31").Value = 35720372
$ws.Range("K31").Value = 1986.7333
$ws.Range("L31").Value = 35720372
$ws.Range("M31").Value = -1691.7333
$ws.Range("N31").Value = -35720962

$ws.Range("H34").Value = 17245344
$ws.Range("I34").Value = 1986.7333
$ws.Range("J34").Value = 35720372
$ws.Range("K34").Value = 1986.7333
$ws.Range("L34").Value = 35720372
$ws.Range("M34").Value = -1784.7333
$ws.Range("N34").Value = -35720776

$ws.Range("H94").Value = 1229.8636
$ws.Range("I94").Value = 828
$ws.Range("J94").Value = 1319.1666
$ws.Range("K94").Value = 828
$ws.Range("L94").Value = 1319.1666
$ws.Range("M94").Value = -377
$ws.Range("N94").Value = -2221.1666

$ws.Range("H107").Value = 790.1875
$ws.Range("I107").Value = 558.63635
$ws.Range("K107").Value = 558.63635
$ws.Range("M107").Value = 1361.36365

$ws.Range("H132").Value = 3796.25
$ws.Range("I132").Value = 4035.5
$ws.Range("K132").Value = 12106.5
$ws.Range("M132").Value = -9576.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3702750
$ws.Range("I11").Value = 7200000
$ws.Range("J11").Value = 205499.9
$ws.Range("K11").Value = 7200000
$ws.Range("L11").Value = 205499.9
$ws.Range("M11").Value = -7199861
$ws.Range("N11").Value = -205777.9

$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H113").Value = 1439.1538
$ws.Range("J113").Value = 1601.8572
$ws.Range("L113").Value = 1601.8572
$ws.Range("N113").Value = -5941.8572

$ws.Range("H122").Value = 2785.2058
$ws.Range("I122").Value = 1684.5
$ws.Range("J122").Value = 4357.643
$ws.Range("K122").Value = 5053.5
$ws.Range("L122").Value = 13072.929
$ws.Range("M122").Value = -2603.5
$ws.Range("N122").Value = -17972.929

$ws.Range("H123").Value = 10947.6
$ws.Range("J123").Value = 10947.6
$ws.Range("L123").Value = 10947.6
$ws.Range("N123").Value = -15847.6

$ws.Range("H126").Value = 2935.25
$ws.Range("I126").Value = 2951.6702
$ws.Range("J126").Value = 2404.3333
$ws.Range("K126").Value = 8855.0106
$ws.Range("L126").Value = 7212.999899999999
$ws.Range("M126").Value = -6385.0106
$ws.Range("N126").Value = -12152.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H17").Value = 4245
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 4993.3335
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 4993.3335
$ws.Range("M17").Value = -1830
$ws.Range("N17").Value = -5333.3335
